$d = $word.ActiveDocument

$replacements = @(
    @{old="134÷9="; new="183÷3="},
    @{old="762÷4="; new="311÷7="},
    @{old="339÷8="; new="595÷7="},
    @{old="607÷6="; new="446÷5="},
    @{old="799÷8="; new="793÷6="},
    @{old="491÷7="; new="592÷5="},
    @{old="884÷7="; new="315÷4="},
    @{old="918÷7="; new="582÷7="},
    @{old="779÷8="; new="849÷2="},
    @{old="813÷3="; new="880÷2="},
    @{old="799÷9="; new="703÷5="},
    @{old="713÷3="; new="753÷8="},
    @{old="360÷9="; new="115÷4="},
    @{old="837÷6="; new="614÷3="},
    @{old="886÷7="; new="821÷7="},
    @{old="584÷5="; new="865÷6="},
    @{old="638÷9="; new="613÷6="},
    @{old="681÷9="; new="254÷3="},
    @{old="806÷8="; new="890÷7="},
    @{old="578÷9="; new="732÷3="},
    @{old="824÷5="; new="731÷7="},
    @{old="844÷9="; new="515÷2="},
    @{old="775÷9="; new="606÷5="},
    @{old="913÷2="; new="714÷8="},
    @{old="661÷7="; new="463÷7="}
)

foreach ($r in $replacements) {
    $range = $d.Content
    $range.Find.Execute($r.old, $true, $false, $false, $false, $false, $true, 1, $false, $r.new, 2)
}
